$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'33.829.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "'1.779.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'224.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("E10").Value = "  -5.28%  "

$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").Value = "'2.037.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("D13").Value = "'11.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.13%  "

$ws.Range("D14").Value = "'1.779.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.52%  "

$ws.Range("D15").Value = "'33.875.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "'0.611"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.11%  "

$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").Value = "'66.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.28%  "

$ws.Range("D19").Value = "'238.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.30%  "

$ws.Range("D20").Value = "'0.0₃0773"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "

$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "'10.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.53%  "

$ws.Range("E23").Value = "  -1.93%  "

$ws.Range("E24").Value = "  -2.34%  "

$ws.Range("D25").Value = "'160.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("D27").Value = "'16.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.96%  "

$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("E31").Value = "  -2.62%  "

$ws.Range("E32").Value = "  -3.60%  "

$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("E34").Value = "  -1.26%  "

$ws.Range("D35").Value = "'1.385.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.05%  "

$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("E37").Value = "  -1.62%  "

$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("E39").Value = "  +2.28%  "

$ws.Range("E40").Value = "  +4.80%  "

$ws.Range("D41").Value = "'78.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "

$ws.Range("D42").Value = "'0.910"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.53%  "

$ws.Range("D43").Value = "'13.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.78%  "

$ws.Range("E44").Value = "  -3.05%  "

$ws.Range("D45").Value = "'0.0₆0140"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.98%  "

$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("E47").Value = "  +3.09%  "

$ws.Range("D48").Value = "'107.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("E49").Value = "  -1.77%  "

$ws.Range("D50").Value = "'1.938.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("E51").Value = "  +0.03%  "
